$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Export template layout adds COLLEGE_CODE / COURSE_CODE columns right
# after S_NO (column A), pushing the existing DUMMY_NUMBER.. columns two
# places to the right.
$ws.Columns("B:C").Insert()

# Header text for the two newly inserted columns
$ws.Range("B1").Value = "COLLEGE_CODE"
$ws.Range("C1").Value = "COURSE_CODE"

# Match the look of the rest of row 1's header cells (bold, centered Arial)
$ws.Range("B1:C1").Font.Name = "Arial"
$ws.Range("B1:C1").Font.Bold = $true
$ws.Range("B1:C1").HorizontalAlignment = -4108
$ws.Range("B1:C1").VerticalAlignment = -4108

# Give the new columns sensible explicit widths
$ws.Columns("B").ColumnWidth = 17
$ws.Columns("C").ColumnWidth = 15.666666666666666
